$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "Library Elements" sheet: recode two Element Type cells
#    G3: "Remote Function Call" -> "BTP App/Extension"
#    G4: "REST" -> "Program"
# ------------------------------------------------------------------
$wsLib = $wb.Worksheets.Item("Library Elements")
$wsLib.Range("G3").Value = "BTP App/Extension"
$wsLib.Range("G4").Value = "Program"

# ------------------------------------------------------------------
# 2) Delete the stale "System Group" author comment on E3
# ------------------------------------------------------------------
$cmt = $wsLib.Range("E3").Comment
if ($cmt -ne $null) {
    $cmt.Delete()
}

# ------------------------------------------------------------------
# 3) "Element Types" sheet: rename header + re-sort / refresh the
#    list of selectable element types.
# ------------------------------------------------------------------
$wsTypes = $wb.Worksheets.Item("Element Types")
$wsTypes.Range("A1").Value = "Element Types"

$elementTypes = @(
    "BTP App/Extension",
    "Classes/Interface",
    "Classic BAdI Implementation",
    "Custom Fiori Application",
    "Enhancement Implementation",
    "Extra Workbench Object",
    "Function Group",
    "Function Module",
    "Package",
    "Program",
    "Table",
    "Transaction"
)
for ($i = 0; $i -lt $elementTypes.Length; $i++) {
    $wsTypes.Cells.Item($i + 2, 1).Value = $elementTypes[$i]
}

# ------------------------------------------------------------------
# 4) "Library Type" sheet: drop the unused, empty third column
# ------------------------------------------------------------------
$wsLibType = $wb.Worksheets.Item("Library Type")
$wsLibType.Columns.Item(3).Delete()

# ------------------------------------------------------------------
# 5) Add a dropdown (list) data validation on the Element Type
#    column of "Library Elements", sourced from "Element Types"
# ------------------------------------------------------------------
$dvRange = $wsLib.Range("G2:G10")
$dvRange.Validation.Add(3, 1, 1, "='Element Types'!`$A`$2:`$A`$13")
$dvRange.Validation.InputTitle = "Choose Element Type"
$dvRange.Validation.InputMessage = " "
$dvRange.Validation.ShowInput = $true
$dvRange.Validation.ShowError = $true
$dvRange.Validation.IgnoreBlank = $true
